# "first architecture check with jochen on friday, 10_06_2016"
# Checkliste_Dokumente.xlsx - add an "Architekturdokument" results column (and
# placeholders for three more document types), fill in the checklist answers,
# tidy up the question wording, and drop the old two-column merged layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# ------------------------------------------------------------------
# 1) Get rid of the old B:C merged "notes" block (rows 17-20) and the
#    centered alignment that came with it, before we rebuild the area.
# ------------------------------------------------------------------
$ws.Range("B17:C17").UnMerge()
$ws.Range("B18:C18").UnMerge()
$ws.Range("B19:C19").UnMerge()
$ws.Range("B20:C20").UnMerge()

$ws.Range("B2").Copy()
$ws.Range("B17:C20").PasteSpecial(-4122)   # xlPasteFormats -> drop centered look
$ws.Application.CutCopyMode = $false

# ------------------------------------------------------------------
# 2) Rewrite the checklist questions in column B (wording fixes, new
#    question added, order changed slightly).
# ------------------------------------------------------------------
$ws.Range("B2").Value  = "Template vom PM benutzt?"
$ws.Range("B3").Value  = "Seitenzahlen vorhanden?"
$ws.Range("B4").Value  = "Seitenzahlen an gleicher Position"
$ws.Range("B5").Value  = "Dokumentenverantwortliche[r] auf Titelseite genannt?"
$ws.Range("B6").Value  = "Titel des Dokuments auf Titelseite?"
$ws.Range("B7").Value  = "Änderungshistorie vorhanden?"
$ws.Range("B8").Value  = "Tabellenverzeichnis vorhanden?"
$ws.Range("B9").Value  = "Abbildungsverzeichnis vorhanden?"
$ws.Range("B10").Value = "Kapitel nummeriert?"
$ws.Range("B11").Value = "Kapitelschachtelung sinnvoll?"
$ws.Range("B12").Value = "Sprache auf Englisch?"
$ws.Range("B13").Value = "Leere Füllseiten im Dokument?"
$ws.Range("B14").Value = "Bilder- und Diagramme lesbar?"
$ws.Range("B15").Value = "Tabellenbeschriftungen vorhanden und sinnvoll?"
$ws.Range("B16").Value = "Abbildungsbeschriftungen vorhanden und sinnvoll?"
$ws.Range("B17").Value = "Zu jedem Diagramm/Abbildung eine Beschreibung vorhanden?"

# Row 19 no longer carries any B/C content at all.
$ws.Range("B19").Clear()
$ws.Range("C19").Clear()

# ------------------------------------------------------------------
# 3) Drop the old (now unused) helper column C content for rows 2-17.
# ------------------------------------------------------------------
$ws.Range("C2:C17").Clear()

# ------------------------------------------------------------------
# 4) Add the new "Architekturdokument" column with the first set of
#    check results, plus the three extra (still empty) document-type
#    headers, and the "checked on" note.
# ------------------------------------------------------------------
$ws.Range("D1").Value = "Architekturdokument"
$ws.Range("E1").Value = "Projekthandbuch"
$ws.Range("F1").Value = "Anforderungsdokument"
$ws.Range("G1").Value = "Testreport(?)"

$ws.Range("D2").Value  = "ja, angepasst"
$ws.Range("D3").Value  = "ja"
$ws.Range("D4").Value  = "ja"
$ws.Range("D5").Value  = "ja"
$ws.Range("D6").Value  = "ja"
$ws.Range("D7").Value  = "ja"
$ws.Range("D8").Value  = "ja"
$ws.Range("D9").Value  = "ja"
$ws.Range("D10").Value = "ja"
$ws.Range("D11").Value = "ja"
$ws.Range("D12").Value = "ja"
$ws.Range("D13").Value = "nein"
$ws.Range("D14").Value = "ja"
$ws.Range("D15").Value = "ja"
$ws.Range("D16").Value = "ja"
$ws.Range("D17").Value = "ja"
$ws.Range("D18").Value = "gecheckt am 10.06.2016"

# D19/E19: blank placeholder cells (matching B18/C18/B20/C20's look).
$ws.Range("D19").Value = "x"
$ws.Range("E19").Value = "x"
$ws.Range("D19:E19").Copy()
$ws.Range("D19:E19").PasteSpecial(-4122)   # re-stamp formats (keeps cell alive)
$ws.Application.CutCopyMode = $false
$ws.Range("D19:E19").ClearContents()

# ------------------------------------------------------------------
# 5) Column widths for the new / widened columns.
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 57.85
$ws.Columns.Item(4).ColumnWidth = 21.85
$ws.Columns.Item(5).ColumnWidth = 16.28
$ws.Columns.Item(6).ColumnWidth = 22.57
$ws.Columns.Item(7).ColumnWidth = 12.71

# ------------------------------------------------------------------
# 6) Leave the selection where the author last left it.
# ------------------------------------------------------------------
$ws.Range("D15").Select()
